# Adds a new "Set RunPeriod" boolean argument row to the AddMonthlyJSONUtilityData
# measure block on the Variables sheet, shifting all the rows below it down by one
# (mirrors inserting a whole worksheet row above the existing row 13).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Variables")

# Insert a new whole row above row 13 - this pushes every row from 13 down by
# one (so the old row 13 "AddMonthlyJSONUtilityDataGas" header becomes row 14,
# etc.) and copies the row-12 formatting into the freshly inserted row 13.
$ws.Rows(13).Insert()

# Populate the new argument row with the same shape as the other "argument"
# rows in this measure block (see row 12 directly above it).
$ws.Range("B13").Value2 = "argument"
$ws.Range("D13").Value2 = "Set RunPeriod"
$ws.Range("E13").Value2 = "set_runperiod"
$ws.Range("G13").Value2 = "Bool"
$ws.Range("I13").Value2 = "TRUE"

# Select the newly inserted row, matching what Excel leaves selected right
# after a row insert.
$ws.Rows(13).Select()

# The sheet already had an AutoFilter over A2:AA126; growing the data by one
# row means the filter (and its hidden _FilterDatabase defined name) needs to
# cover A2:AA127 instead. Toggling AutoFilterMode off/on lets us re-apply it
# over the new range.
$ws.AutoFilterMode = $false
$ws.Range("A2:AA127").AutoFilter()

foreach ($n in $wb.Names) {
    if ($n.Name -eq "Variables!_FilterDatabase") {
        $n.RefersTo = "=Variables!`$A`$2:`$AA`$127"
    }
}
